# Scheduled runner: refresh cached market-board price/profit snapshots
# across each crafting-class sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) in the
# Leve profit workbook. Values below are the newly-fetched figures; where a
# source no longer reports a figure the cell is cleared instead of zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 34.285713
$ws.Range("I53").Value = 42.75
$ws.Range("J53").Value = 23
$ws.Range("K53").Value = 42.75
$ws.Range("L53").Value = 23
$ws.Range("M53").Value = 594.25
$ws.Range("N53").Value = -1297

# Row 135
$ws.Range("H135").Value = 1461.4286
$ws.Range("I135").Value = 1371.6666
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 12344.9994
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -9809.999400000001
$ws.Range("N135").Value = -23070

# Row 137
$ws.Range("H137").Value = 1390.4783
$ws.Range("J137").Value = 1580.2
$ws.Range("L137").Value = 4740.6
$ws.Range("N137").Value = -9840.6

# Row 138
$ws.Range("H138").Value = 2066.4531
$ws.Range("I138").Value = 1815
$ws.Range("J138").Value = 2180.75
$ws.Range("K138").Value = 5445
$ws.Range("L138").Value = 6542.25
$ws.Range("M138").Value = -305
$ws.Range("N138").Value = -16822.25
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11825.652
$ws.Range("I32").Value = 9111
$ws.Range("K32").Value = 9111
$ws.Range("M32").Value = -8824

# Row 45
$ws.Range("H45").Value = 3237.375
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 3414.2856
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 3414.2856
$ws.Range("M45").Value = -1622
$ws.Range("N45").Value = -4168.2856

# Row 132
$ws.Range("H132").Value = 4249.4546
$ws.Range("I132").Value = 2999.25
$ws.Range("K132").Value = 8997.75
$ws.Range("M132").Value = -6467.75
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3665
$ws.Range("I20").Value = 3665
$ws.Range("K20").Value = 3665
$ws.Range("M20").Value = -3418

# Row 99
$ws.Range("H99").Value = 1999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1999
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4995

# Row 105
$ws.Range("H105").Value = 3585.739
$ws.Range("I105").Value = 3581.8333
$ws.Range("J105").Value = 3590
$ws.Range("K105").Value = 3581.8333
$ws.Range("L105").Value = 3590
$ws.Range("M105").Value = -1834.8333
$ws.Range("N105").Value = -7084

# Row 107
$ws.Range("H107").Value = 4233.3335
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4233.3335
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4233.3335
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8073.3335

# Row 134
$ws.Range("H134").Value = 4001.4285
$ws.Range("I134").Value = 3602.2
$ws.Range("K134").Value = 10806.6
$ws.Range("M134").Value = -8271.599999999999
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 25163.5
$ws.Range("I3").Value = 22666.666
$ws.Range("J3").Value = 27660.334
$ws.Range("K3").Value = 22666.666
$ws.Range("L3").Value = 27660.334
$ws.Range("M3").Value = -22553.666
$ws.Range("N3").Value = -27886.334

# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 74
$ws.Range("H74").Value = 38593
$ws.Range("J74").Value = 38593
$ws.Range("L74").Value = 38593
$ws.Range("N74").Value = -40341

# Row 77
$ws.Range("H77").Value = 38593
$ws.Range("J77").Value = 38593
$ws.Range("L77").Value = 115779
$ws.Range("N77").Value = -124515

# Row 105
$ws.Range("H105").Value = 3299.6667
$ws.Range("I105").Value = 2998
$ws.Range("J105").Value = 3360
$ws.Range("K105").Value = 2998
$ws.Range("L105").Value = 3360
$ws.Range("M105").Value = -1251
$ws.Range("N105").Value = -6854

# Row 107
$ws.Range("H107").Value = 2414.6177
$ws.Range("I107").Value = 2991.611
$ws.Range("J107").Value = 1765.5
$ws.Range("K107").Value = 2991.611
$ws.Range("L107").Value = 1765.5
$ws.Range("M107").Value = -1071.611
$ws.Range("N107").Value = -5605.5
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# Row 23
$ws.Range("H23").Value = 200026.8
$ws.Range("J23").Value = 500037
$ws.Range("L23").Value = 1500111
$ws.Range("N23").Value = -1500581

# Row 37
$ws.Range("H37").Value = 99250
$ws.Range("J37").Value = 99250
$ws.Range("L37").Value = 297750
$ws.Range("N37").Value = -297974

# Row 64
$ws.Range("H64").Value = 16337.333
$ws.Range("J64").Value = 14500
$ws.Range("L64").Value = 43500
$ws.Range("N64").Value = -44040

# Row 67
$ws.Range("H67").Value = 16337.333
$ws.Range("J67").Value = 14500
$ws.Range("L67").Value = 43500
$ws.Range("N67").Value = -45372

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 5254
$ws.Range("I4").Value = 6504.5
$ws.Range("J4").Value = 4003.5
$ws.Range("K4").Value = 6504.5
$ws.Range("L4").Value = 4003.5
$ws.Range("M4").Value = -6391.5
$ws.Range("N4").Value = -4229.5

# Row 28
$ws.Range("H28").Value = 5254
$ws.Range("I28").Value = 6504.5
$ws.Range("J28").Value = 4003.5
$ws.Range("K28").Value = 6504.5
$ws.Range("L28").Value = 4003.5
$ws.Range("M28").Value = -6272.5
$ws.Range("N28").Value = -4467.5

# Row 35
$ws.Range("H35").Value = 195
$ws.Range("I35").Value = 195
$ws.Range("K35").Value = 195
$ws.Range("M35").Value = 141

# Row 37
$ws.Range("H37").Value = 5254
$ws.Range("I37").Value = 6504.5
$ws.Range("J37").Value = 4003.5
$ws.Range("K37").Value = 6504.5
$ws.Range("L37").Value = 4003.5
$ws.Range("M37").Value = -6397.5
$ws.Range("N37").Value = -4217.5

# Row 122
$ws.Range("H122").Value = 3086.8
$ws.Range("J122").Value = 2922
$ws.Range("L122").Value = 8766
$ws.Range("N122").Value = -13666

# Row 136
$ws.Range("H136").Value = 6087.6
$ws.Range("I136").Value = 5748.75
$ws.Range("J136").Value = 7443
$ws.Range("K136").Value = 17246.25
$ws.Range("L136").Value = 22329
$ws.Range("M136").Value = -14696.25
$ws.Range("N136").Value = -27429
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 378.75
$ws.Range("I113").Value = 305.22223
$ws.Range("K113").Value = 915.66669
$ws.Range("M113").Value = 1254.33331

# Row 136
$ws.Range("H136").Value = 1467.2858
$ws.Range("I136").Value = 776.7273
$ws.Range("K136").Value = 2330.1819
$ws.Range("M136").Value = 219.8181
